# Apply cryptos list update (values refreshed on Tue Nov 19 19:53:05 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.454.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.130.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "617.06"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.408"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +10.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.128.13"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +30.95%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.72%  "

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.193.53"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.14%  "

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.85"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.727.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.074.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.81"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.93"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.86"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000208"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "450.90"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.63%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.87"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.91"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.302.43"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +13.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.231"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.30"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.88%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.43"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "493.22"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.86%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.18%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.437"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.47"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.11"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.53%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.94"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.01%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.84%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0336"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.48"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.24%  "
